# Remove the "link tik tok" and "link github" paragraphs (each contains a
# hyperlink run) that were added at the end of the document, restoring the
# trailing structure to: <empty noProof paragraph> <empty paragraph> <sectPr>.

$d = $word.ActiveDocument

# Walk the paragraphs from the end towards the start so deleting one
# paragraph does not disturb the indices of paragraphs we still need to
# inspect.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    $isTikTokLink = $text -like "link tik tok:*tiktok.com*"
    $isGithubLink = $text -like "link*github*github.com*"

    if ($isTikTokLink -or $isGithubLink) {
        $para.Range.Delete()
    }
}
